# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right after "总计" (so order becomes
#   总计, 2022-Q3, 2022-Q2, 2022-Q1) and populate it with the quarterly
#   fund-holding detail rows.
# - Update the "总计" (summary) sheet so it now lists three rows
#   (2022-Q3, 2022-Q2, 2022-Q1) instead of two.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift the existing two data rows down by one, then write in the new
# 2022-Q3 row at the top of the data block.
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(4, 1).PasteSpecial(-4122)
$summary.Cells.Item(4, 2).Value = "2022-Q1"
$summary.Cells.Item(4, 3).Value = 1
$summary.Cells.Item(4, 4).Value = 0.03

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 3
$summary.Cells.Item(3, 4).Value = 0.06

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 7
$summary.Cells.Item(2, 4).Value = 0.18

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# "2022-Q2" (now 3rd tab, right after the freshly-inserted "2022-Q3") is
# the formatting template for the new quarter sheet - identical
# header/column layout.
$template = $wb.Worksheets.Item(3)

# Header row - copy formatting (bold/border/center) from the template
# sheet's header, then fill in the text.
$template.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Fund-holding detail rows. Columns D-G are kept as text (matching the
# other quarter sheets), so force a text number format before writing
# then drop back to the default cell style afterwards.
$q3Data = @(
    @(0, "470888", "汇添富香港优势精选混合（QDII）", "1.63", "78.50", "6.83", "0.1113", 3),
    @(1, "007254", "广发均衡价值混合",               "0.66", "89.74", "7.86", "0.0519", 1),
    @(2, "012315", "创金合信港股通成长股票A",         "0.08", "80.48", "8.78", "0.0070", 5),
    @(3, "012316", "创金合信港股通成长股票C",         "0.07", "80.48", "8.78", "0.0061", 5),
    @(4, "006603", "嘉实互融精选股票",                "0.12", "82.85", "2.61", "0.0031", 9),
    @(5, "013182", "安信港股通精选混合C",             "0.12", "38.51", "1.26", "0.0015", 10),
    @(6, "013181", "安信港股通精选混合A",             "0.02", "38.51", "1.26", "0.0003", 10)
)

$rowNum = 2
foreach ($rec in $q3Data) {
    # Columns B-G hold text values in the source data (fund code keeps its
    # leading zeros, percentages/amounts are plain text) - force a text
    # number format before writing, then drop back to the default style.
    $rng = $q3.Range($q3.Cells.Item($rowNum, 2), $q3.Cells.Item($rowNum, 7))
    $rng.NumberFormat = "@"

    $q3.Cells.Item($rowNum, 1).Value = $rec[0]
    $q3.Cells.Item($rowNum, 2).Value = $rec[1]
    $q3.Cells.Item($rowNum, 3).Value = $rec[2]
    $q3.Cells.Item($rowNum, 4).Value = $rec[3]
    $q3.Cells.Item($rowNum, 5).Value = $rec[4]
    $q3.Cells.Item($rowNum, 6).Value = $rec[5]
    $q3.Cells.Item($rowNum, 7).Value = $rec[6]
    $q3.Cells.Item($rowNum, 8).Value = $rec[7]

    $rng.Style = "Normal"

    $template.Cells.Item(2, 1).Copy()
    $q3.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    $q3.Cells.Item($rowNum, 1).Value = $rec[0]

    $rowNum = $rowNum + 1
}

# ---------------------------------------------------------------------
# 3. Restore the originally-active sheet ("2022-Q1"), since inserting a
#    new worksheet shifts focus onto it.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
